$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..18 (each row r takes A/B/C/D from old row r+1, with C
# recomputed, and a freshly recalculated E value; old row 19 is dropped).
$data = @(
    @(2, 39765, 2008, $null, 2009, $null),
    @(3, 40130, 2009, 1.834695583582535, 2010, $null),
    @(4, 40494, 2010, 1.767835936772144, 2011, $null),
    @(5, 40862, 2011, 1.074400434091038, 2012, $null),
    @(6, 41228, 2012, 0.9212998022035679, 2013, 1.274704633957136),
    @(7, 41592, 2013, 1.141837882844188, 2014, 1.404348988410131),
    @(8, 41957, 2014, 1.335361538769475, 2015, 1.269653854937691),
    @(9, 42321, 2015, 1.202048372526998, 2016, 1.253742200752095),
    @(10, 42689, 2016, 2.677488680362305, 2017, 1.805615391969595),
    @(11, 43053, 2017, 2.466954516646402, 2018, 1.661541796722577),
    @(12, 43418, 2018, 1.401189216021326, 2019, 1.815016201748643),
    @(13, 43783, 2019, 2.217567799050979, 2020, 1.810449264563152),
    @(14, 44159, 2020, 2.139672475020404, 2021, 2.128328071999674),
    @(15, 44525, 2021, 2.100991693542231, 2022, 1.11435041103376),
    @(16, 44890, 2022, 0.8967077601845341, 2023, 0.3338851812143995),
    @(17, 45254, 2023, 0.782207885866093, 2024, 2.228542839642689),
    @(18, 45618, 2024, 1.508385007449875, 2025, 0.9823016603409229)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    if ($row[3] -ne $null) {
        $ws.Cells.Item($r, 3).Value = $row[3]
    } else {
        $ws.Cells.Item($r, 3).Value = $null
    }
    $ws.Cells.Item($r, 4).Value = $row[4]
    if ($row[5] -ne $null) {
        $ws.Cells.Item($r, 5).Value = $row[5]
    } else {
        $ws.Cells.Item($r, 5).Value = $null
    }
}

# Remove the now-unused last row (old row 19) entirely.
$ws.Rows.Item(19).Delete()
